{"js": "// Add a hyperlink paragraph (linking to the explanatory video) right after\n// the \"Link para o v\u00eddeo explicativo:\" paragraph, followed by a blank\n// paragraph \u2014 mirroring the existing \"Link para a simula\u00e7\u00e3o\"/\"Link para o\n// Projeto Github\" hyperlink paragraphs already in the document.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the paragraph that introduces the video link.\nconst items = paragraphs.items;\nlet anchorParagraph = items[items.length - 1];\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"Link para o v\\u00eddeo explicativo:\") !== -1) {\n    anchorParagraph = items[i];\n    break;\n  }\n}\n\nconst videoUrl = \"https://youtu.be/gy_Ql1mnMQY\";\n\n// New paragraph holding the hyperlink text.\nconst linkParagraph = anchorParagraph.insertParagraph(videoUrl, Word.InsertLocation.after);\nconst linkRange = linkParagraph.getRange();\nlinkRange.hyperlink = videoUrl;\n\n// Trailing blank paragraph, same as after the other hyperlinks in the doc.\nlinkParagraph.insertParagraph(\"\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Add a hyperlink paragraph (linking to the explanatory video) right after\n# the \"Link para o v\u00eddeo explicativo:\" paragraph, followed by a blank\n# paragraph \u2014 mirroring the existing \"Link para a simula\u00e7\u00e3o\"/\"Link para o\n# Projeto Github\" hyperlink paragraphs already in the document.\n\n$d = $word.ActiveDocument\n$videoUrl = \"https://youtu.be/gy_Ql1mnMQY\"\n\n# Locate the paragraph that introduces the video link.\n$found = $d.Content.Find.Execute(\"Link para o v\u00eddeo explicativo:\")\nif ($found) {\n    $anchorParaIndex = $d.Paragraphs.Count\n    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n        if ($d.Paragraphs.Item($i).Range.Text -like \"*Link para o v\u00eddeo explicativo:*\") {\n            $anchorParaIndex = $i\n        }\n    }\n} else {\n    $anchorParaIndex = $d.Paragraphs.Count\n}\n\n# Insert a new (initially empty) paragraph right after the anchor paragraph.\n$anchorPara = $d.Paragraphs.Item($anchorParaIndex)\n$insPoint = $anchorPara.Range\n$insPoint.Collapse(0) | Out-Null   # wdCollapseEnd\n$insPoint.InsertParagraphAfter() | Out-Null\n\n# The just-inserted paragraph is now the last one; its Start is the caret\n# position inside it (re-read live from the Paragraphs collection \u2014 a Range\n# object held across a mutating call goes stale).\n$linkParaIndex = $d.Paragraphs.Count\n$startPos = $d.Paragraphs.Item($linkParaIndex).Range.Start\n\n$linkRange = $d.Range($startPos, $startPos)\n$linkRange.InsertAfter($videoUrl) | Out-Null\n\n# Turn the freshly typed URL text into a real hyperlink.\n$endPos = $startPos + $videoUrl.Length\n$hlTarget = $d.Range($startPos, $endPos)\n$d.Hyperlinks.Add($hlTarget, $videoUrl) | Out-Null\n\n# Trailing blank paragraph, same as after the other hyperlinks in the doc.\n$tailParaIndex = $d.Paragraphs.Count\n$tailPoint = $d.Paragraphs.Item($tailParaIndex).Range\n$tailPoint.Collapse(0) | Out-Null\n$tailPoint.InsertParagraphAfter() | Out-Null\n"}
